$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-07 Friday" "2025-11-08 Saturday"

Replace-Text "83×47=" "93×95="
Replace-Text "79×38=" "36×35="
Replace-Text "70×77=" "62×46="
Replace-Text "76×74=" "91×33="
Replace-Text "11×49=" "88×55="

Replace-Text "23×42=" "47×32="
Replace-Text "65×19=" "88×85="
Replace-Text "27×58=" "42×93="
Replace-Text "75×56=" "40×96="
Replace-Text "92×76=" "73×11="

Replace-Text "53×92=" "59×51="
Replace-Text "75×35=" "20×83="
Replace-Text "44×26=" "40×75="
Replace-Text "35×30=" "24×82="
Replace-Text "18×74=" "39×25="

Replace-Text "88×81=" "37×33="
Replace-Text "98×79=" "84×71="
Replace-Text "23×49=" "93×34="
Replace-Text "96×41=" "40×21="
Replace-Text "67×22=" "63×97="

Replace-Text "91×62=" "85×70="
Replace-Text "16×53=" "89×47="
Replace-Text "62×54=" "17×43="
Replace-Text "30×40=" "55×73="
Replace-Text "16×18=" "15×83="
